# Apply the coin-price / volume refresh captured in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.988.58"
$ws.Range("E2").Value = "  +0.19%  "

$ws.Range("D3").Value = "3.422.73"
$ws.Range("E3").Value = "  +0.34%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").Value = "'409.90"
$ws.Range("E5").Value = "  +1.41%  "

$ws.Range("D6").Value = "'128.52"
$ws.Range("E6").Value = "  -3.42%  "

$ws.Range("E7").Value = "  +7.37%  "

$ws.Range("D8").Value = "'0.999"
$ws.Range("E8").Value = "  -0.10%  "

$ws.Range("E9").Value = "  +6.77%  "

$ws.Range("E10").Value = "  +10.99%  "

$ws.Range("D11").Value = "'42.60"
$ws.Range("E11").Value = "  +1.63%  "

$ws.Range("E12").Value = "  +0.17%  "

$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'9.04"
$ws.Range("E13").Value = "  +8.01%  "

$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "3.955.71"
$ws.Range("E14").Value = "  -0.06%  "

$ws.Range("D15").Value = "'21.15"
$ws.Range("E15").Value = "  +7.06%  "

$ws.Range("E16").Value = "  +45.96%  "

$ws.Range("D17").Value = "3.410.54"
$ws.Range("E17").Value = "  -0.47%  "

$ws.Range("D18").Value = "'12.37"
$ws.Range("E18").Value = "  +5.71%  "

$ws.Range("E19").Value = "  +6.58%  "

$ws.Range("D20").Value = "61.905.50"
$ws.Range("E20").Value = "  -0.06%  "

$ws.Range("D21").Value = "'447.21"
$ws.Range("E21").Value = "  +43.79%  "

$ws.Range("D22").Value = "'91.72"
$ws.Range("E22").Value = "  +10.32%  "

$ws.Range("D23").Value = "'3.18"
$ws.Range("E23").Value = "  +0.79%  "

$ws.Range("D24").Value = "'12.90"
$ws.Range("E24").Value = "  +1.53%  "

$ws.Range("E25").Value = "  +2.90%  "

$ws.Range("D26").Value = "'33.14"
$ws.Range("E26").Value = "  +12.08%  "

$ws.Range("E27").Value = "  +9.27%  "

$ws.Range("D28").Value = "'4.79"
$ws.Range("E28").Value = "  -0.56%  "

$ws.Range("E29").Value = "  -0.81%  "

$ws.Range("E30").Value = "  -5.37%  "

$ws.Range("D31").Value = "'11.96"
$ws.Range("E31").Value = "  +5.40%  "

$ws.Range("D32").Value = "'0.170"
$ws.Range("E32").Value = "  -1.37%  "

$ws.Range("D33").Value = "'0.114"
$ws.Range("E33").Value = "  -0.95%  "

$ws.Range("D34").Value = "'42.74"
$ws.Range("E34").Value = "  +0.69%  "

$ws.Range("E35").Value = "  -0.06%  "

$ws.Range("D36").Value = "'0.0497"
$ws.Range("E36").Value = "  +3.21%  "

$ws.Range("D37").Value = "'53.28"
$ws.Range("E37").Value = "  +3.75%  "

$ws.Range("D38").Value = "'0.998"
$ws.Range("E38").Value = "  -0.10%  "

$ws.Range("D39").Value = "'3.38"
$ws.Range("E39").Value = "  +0.40%  "

$ws.Range("D40").Value = "'0.134"
$ws.Range("E40").Value = "  +7.72%  "

$ws.Range("D41").Value = "'2.92"
$ws.Range("E41").Value = "  -1.24%  "

$ws.Range("D42").Value = "'0.314"
$ws.Range("E42").Value = "  -2.39%  "

$ws.Range("D43").Value = "'141.13"
$ws.Range("E43").Value = "  +0.86%  "

$ws.Range("D44").Value = "'4.24"
$ws.Range("E44").Value = "  +7.49%  "

$ws.Range("D45").Value = "'1.98"
$ws.Range("E45").Value = "  +0.67%  "

$ws.Range("E46").Value = "  +8.23%  "

$ws.Range("D47").Value = "'16.53"
$ws.Range("E47").Value = "  -0.29%  "

$ws.Range("D48").Value = "'22.43"
$ws.Range("E48").Value = "  +5.25%  "

$ws.Range("D49").Value = "3.765.65"
$ws.Range("E49").Value = "  +0.36%  "

$ws.Range("D50").Value = "'2.09"
$ws.Range("E50").Value = "  +8.15%  "

$ws.Range("D51").Value = "2.118.88"
$ws.Range("E51").Value = "  +0.75%  "

Write-Output "Applied cryptos update to 50 rows"